$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update row 26 (2026-01) stats
$ws.Range("B26").Value = 6531
$ws.Range("C26").Value = 1012
$ws.Range("D26").Value = 6079887
$ws.Range("E26").Value = 930.9274230592558
$ws.Range("F26").Value = 10.37688017576475
$ws.Range("G26").Value = 7.430997876857748
$ws.Range("H26").Value = 26.61385160455843
